# daily auto push: 2026-01-10 22:32 UTC
# A new timestamped observation (2026/01/11, 日, hour 5, ranking 19) was
# appended right after the existing 2026/01/11 block, which pushes every
# subsequent row (the old row 599 .. row 640) down by one. The sheet's
# used range grows from A1:D640 to A1:D641.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 599..640 down to 600..641, opening up a blank row at 599.
$ws.Rows.Item(599).Insert()

# Fill in the new row. Column A holds a date-like string ("2026/01/11")
# that must stay literal text (matching the rest of the column, which is
# plain text, not a real date serial) -- a leading apostrophe forces that
# the same way it does when typing directly into Excel.
$ws.Cells.Item(599, 1).Value = "'2026/01/11"
$ws.Cells.Item(599, 2).Value = "日"
$ws.Cells.Item(599, 3).Value = 5
$ws.Cells.Item(599, 4).Value = 19
